# Daily update of covid19 tracker data files
# - Bumps the "Updated on" date column (B5:B74) from 6-Apr-2020 to 7-Apr-2020
# - Applies small copy edits to a handful of country narrative cells

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Updates")

# --- 1. Advance the "Updated on" date by one day for every data row (5-74) ---
for ($r = 5; $r -le 74; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = $cell.Value2 + 1
}

# --- 2. Czech Republic (J16): remove the stray double space before "launched" ---
$old = $ws.Range("J16").Value2
$new = $old.Replace("(Antivirus,  launched April 6)", "(Antivirus, launched April 6)")
$ws.Range("J16").Value = $new

# --- 3. Germany (J22): "10 - 50 employees" -> "10 to 50 employees" ---
$old = $ws.Range("J22").Value2
$new = $old.Replace("EUR 800 000 for firms with 10 - 50 employees", "EUR 800 000 for firms with 10 to 50 employees")
$ws.Range("J22").Value = $new

# --- 4. India (G26): clarify the hydroxychloroquine export exemption wording ---
$old = $ws.Range("G26").Value2
$new = $old.Replace("The government withdrew exemptions on exports of hydroxychloroquine", "The government withdrew exemptions from the ban on exports of hydroxychloroquine")
$ws.Range("G26").Value = $new

# --- 5. Israel (E29): full rewrite of the schools-closure sentence ---
$israelText = @"
Since 13 March all educational institutions including pre-schools and kindergartens are closed at least until the end of Passover (17 April)
"@
$ws.Range("E29").Value = $israelText

# --- 6. Turkey (F53): drop the Istanbul-governor sentence, add mask distribution info ---
$turkeyText = @"
Persons above 65 (23 March) and those below 20 (4 April) are locked-down in their living places (with the exception of young workers between 18-20). Sport and cultural facilities, mosques (during high attendance periods), domestic and international trade fairs, cafes, museums, librarie.  are closed.  From 28 March, outside physical exercises and picnics are banned in town centers on the weekends. Local authorities can extend these bans to weekdays. Grocery stores and supermarkets' opening hours are restricted, with a maximum of one customer for every 10 square metres of  space. From 4 April, wearing a protective masks is made compulsory in all public places. On 6 April the government announced that it will distribute, to all citizens between 20-65 placing a request via the e-government portal, 5 free protective masks per week via postal delivery.
"@
$ws.Range("F53").Value = $turkeyText
